$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the post row for "「シベリアジャコウジカは、絶滅の危機にあると見なされている」"
# (row 420). All rows below it shift up by one, and the used range shrinks
# from A1:C472 to A1:C471.
$ws.Rows.Item(420).Delete()
